$d = $word.ActiveDocument

# --- 1. Update the subtitle date field text ---
$d.Content.Find.Execute("09 October 2015", $false, $false, $false, $false, $false, $true, 1, $false, "12 October 2015", 2) | Out-Null

# --- 2. Move the "_GoBack" bookmark from the title paragraph down to just before
#        the final ". " run near the end of the References section ---
$rEnd = $d.Content
$rEnd.Find.Execute("www.GenStat.co.uk", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackTarget = $d.Range($rEnd.End, $rEnd.End)
$d.Bookmarks.Add("_GoBack", $goBackTarget) | Out-Null

# --- 3. Rewrite the "To install R..." / "Important:" paragraphs into a
#        "Follow the steps below..." intro + Step 1-4 paragraphs ---
$rA = $d.Content
$rA.Find.Execute("To install R, go to the R website", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rB = $d.Content
$rB.Find.Execute('Save version number in registry" checked.', $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$installRRange = $d.Range($rA.Start, $rB.End)
$installRXml = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>    <w:p>
      <w:r>
        <w:t>Follow the steps below t</w:t>
      </w:r>
      <w:r>
        <w:t>o install R</w:t>
      </w:r>
      <w:r>
        <w:t>:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Step 1:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>G</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">o to the R website for downloading the Windows version on </w:t>
      </w:r>
      <w:hyperlink r:id="rId11" w:history="1">
        <w:r>
          <w:t>http://cran.rstudio.org</w:t>
        </w:r>
      </w:hyperlink>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Step 2:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Click on the link "Download </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>R.x.x.x</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> for Windows" (or other version). This starts downloading R.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>x</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>x</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>x</w:t>
      </w:r>
      <w:r>
        <w:t>-win.e</w:t>
      </w:r>
      <w:r>
        <w:t>xe file for both 32 and 64 bit.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Step 3:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>After downloading, double click this file to ins</w:t>
      </w:r>
      <w:r>
        <w:t>tall R.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Important:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Make sure that you keep the default setting under Additional Tasks: "Save version number in registry" checked.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:i/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Step 4:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>S</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">tart R and install the following packages </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>MASS</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>lsmeans</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>stringr</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>reshape</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>mvtnorm</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> by typing:  </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>install.packages</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>("</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>package_name_here</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:i/>
        </w:rPr>
        <w:t>")</w:t>
      </w:r>
      <w:r>
        <w:t>. Thes</w:t>
      </w:r>
      <w:r>
        <w:t>e packages are used for the analysis.</w:t>
      </w:r>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$installRRange.InsertXML($installRXml)

# InsertXML across multiple new paragraphs loses the inner <w:rPr> of a
# w:hyperlink run that is not the first run of the first paragraph, so the
# hyperlink was inserted without its "Hyperlink" character style - fix it up
# with a follow-up Find + Style assignment (single-run range, not affected).
$fixLink = $d.Content
$fixLink.Find.Execute("http://cran.rstudio.org", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fixLink.Style = "Hyperlink"
